$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 315.17392
$ws.Range("I28").Value = 266.1579
$ws.Range("J28").Value = 548
$ws.Range("K28").Value = 266.1579
$ws.Range("L28").Value = 548
$ws.Range("M28").Value = 218.8421
$ws.Range("N28").Value = -1518

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 168.85715
$ws.Range("I33").Value = 143.38461
$ws.Range("K33").Value = 143.38461
$ws.Range("M33").Value = 85.61538999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5992.5
$ws.Range("I43").Value = 5992.5
$ws.Range("K43").Value = 5992.5
$ws.Range("M43").Value = -5923.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3785.6667
$ws.Range("I64").Value = 3608.182
$ws.Range("K64").Value = 3608.182
$ws.Range("M64").Value = -3360.182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3785.6667
$ws.Range("I67").Value = 3608.182
$ws.Range("K67").Value = 3608.182
$ws.Range("M67").Value = -2750.182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1438
$ws.Range("I118").Value = 251.25
$ws.Range("J118").Value = 2624.75
$ws.Range("K118").Value = 753.75
$ws.Range("L118").Value = 7874.25
$ws.Range("M118").Value = 903.25
$ws.Range("N118").Value = -11188.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 6904
$ws.Range("J121").Value = 6904
$ws.Range("L121").Value = 20712
$ws.Range("N121").Value = -24206

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 16422
$ws.Range("I132").Value = 16422
$ws.Range("K132").Value = 49266
$ws.Range("M132").Value = -46736

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4194.9565
$ws.Range("I138").Value = 7458.2
$ws.Range("J138").Value = 3288.5
$ws.Range("K138").Value = 22374.6
$ws.Range("L138").Value = 9865.5
$ws.Range("M138").Value = -17234.6
$ws.Range("N138").Value = -20145.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1144.1666
$ws.Range("I41").Value = 1144.1666
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1144.1666
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -730.1666
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3599.5
$ws.Range("J63").Value = 6899.5
$ws.Range("L63").Value = 6899.5
$ws.Range("N63").Value = -8271.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3599.5
$ws.Range("J66").Value = 6899.5
$ws.Range("L66").Value = 34497.5
$ws.Range("N66").Value = -41361.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 30981.223
$ws.Range("J82").Value = 43549.6
$ws.Range("L82").Value = 43549.6
$ws.Range("N82").Value = -44315.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 30981.223
$ws.Range("J85").Value = 43549.6
$ws.Range("L85").Value = 43549.6
$ws.Range("N85").Value = -46201.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2870.5293
$ws.Range("I86").Value = 3061.8696
$ws.Range("J86").Value = 2470.4546
$ws.Range("K86").Value = 3061.8696
$ws.Range("L86").Value = 2470.4546
$ws.Range("M86").Value = -1938.8696
$ws.Range("N86").Value = -4716.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2870.5293
$ws.Range("I89").Value = 3061.8696
$ws.Range("J89").Value = 2470.4546
$ws.Range("K89").Value = 15309.348
$ws.Range("L89").Value = 12352.273
$ws.Range("M89").Value = -9693.348
$ws.Range("N89").Value = -23584.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4758.1763
$ws.Range("J94").Value = 3375
$ws.Range("L94").Value = 3375
$ws.Range("N94").Value = -4277

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 21700
$ws.Range("I99").Value = 21700
$ws.Range("K99").Value = 21700
$ws.Range("M99").Value = -20202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2677.7354
$ws.Range("I134").Value = 2536.6843
$ws.Range("K134").Value = 7610.0529
$ws.Range("M134").Value = -5075.0529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 233
$ws.Range("I13").Value = 56.666668
$ws.Range("J13").Value = 497.5
$ws.Range("K13").Value = 170.000004
$ws.Range("L13").Value = 1492.5
$ws.Range("M13").Value = -2.00000399999999
$ws.Range("N13").Value = -1828.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 6977.909
$ws.Range("I109").Value = 585.6667
$ws.Range("J109").Value = 9375
$ws.Range("K109").Value = 1757.0001
$ws.Range("L109").Value = 28125
$ws.Range("M109").Value = -717.0001
$ws.Range("N109").Value = -30205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 34.875
$ws.Range("K2").Value = 34.875
$ws.Range("M2").Value = 78.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24395.818
$ws.Range("I70").Value = 31429.375
$ws.Range("J70").Value = 5639.6665
$ws.Range("K70").Value = 31429.375
$ws.Range("L70").Value = 5639.6665
$ws.Range("M70").Value = -31159.375
$ws.Range("N70").Value = -6179.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 24395.818
$ws.Range("I73").Value = 31429.375
$ws.Range("J73").Value = 5639.6665
$ws.Range("K73").Value = 31429.375
$ws.Range("L73").Value = 5639.6665
$ws.Range("M73").Value = -30493.375
$ws.Range("N73").Value = -7511.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2850.5833
$ws.Range("I80").Value = 2713.25
$ws.Range("J80").Value = 3125.25
$ws.Range("K80").Value = 2713.25
$ws.Range("L80").Value = 3125.25
$ws.Range("M80").Value = -1715.25
$ws.Range("N80").Value = -5121.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2850.5833
$ws.Range("I83").Value = 2713.25
$ws.Range("J83").Value = 3125.25
$ws.Range("K83").Value = 13566.25
$ws.Range("L83").Value = 15626.25
$ws.Range("M83").Value = -8574.25
$ws.Range("N83").Value = -25610.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 15000
$ws.Range("K23").Value = 15000
$ws.Range("M23").Value = -14770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3880.75
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 4737.5293
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 4737.5293
$ws.Range("M46").Value = -1612
$ws.Range("N46").Value = -5113.5293

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
